# Apply the two changes captured by the commit:
#   1. Slide 6's table switches to a different (built-in) table style.
#   2. The deck's theme colour scheme is swapped from the "Integral"
#      palette to the standard "Office Theme" palette (dk2, lt2,
#      accent1-6, hlink and folHlink all change; dk1/lt1 stay
#      000000/FFFFFF in both palettes).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Table style on slide 6 (the "SOURCES OF FINANCE" table).
# ---------------------------------------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{601E3A56-298F-491D-A526-F623B773F1ED}")
    }
}

# ---------------------------------------------------------------
# 2) Theme colour scheme: Integral -> Office Theme.
#    Item order (matches msoThemeColorSchemeIndex): dk1, lt1, dk2,
#    lt2, accent1..accent6, hlink, folHlink.
# ---------------------------------------------------------------
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
